$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date for every existing data row (2-536)
# from 45202 (2023-10-03) to 45203 (2023-10-04).
$ws.Range("C2:C536").Value = 45203

# The previously-last row (536) gets an explicit row height written out
# (matches the rest of the data rows which already carry ht="15" customHeight="1").
$ws.Rows.Item(536).RowHeight = 15

# Add the new row 537 with the new cutting notification entry.
$ws.Range("A537").Value = "A 47215-2023"
$ws.Range("B537").Value = 45202
$ws.Range("C537").Value = 45203
$ws.Range("D537").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E537").Value = "MOTALA"
$ws.Range("G537").Value = 5.7
$ws.Range("H537").Value = 0
$ws.Range("I537").Value = 0
$ws.Range("J537").Value = 0
$ws.Range("K537").Value = 0
$ws.Range("L537").Value = 0
$ws.Range("M537").Value = 0
$ws.Range("N537").Value = 0
$ws.Range("O537").Value = 0
$ws.Range("P537").Value = 0
$ws.Range("Q537").Value = 0

# Match date formatting used by the rest of column B/C and the wrap-text
# alignment used by column R for the new row.
$ws.Range("B537:C537").NumberFormat = "YYYY-MM-DD"
$ws.Range("R537").WrapText = $true
